# The document currently has no word/styles.xml part at all (an
# incomplete/hand-trimmed test fixture). Word always expects a Styles
# collection to exist, so mint the part back by adding the default
# "Normal" paragraph style through the Styles collection, exactly like
# Word itself would do the first time a style is defined on a document
# that is missing its styles part.

$d = $word.ActiveDocument

# wdStyleTypeParagraph = 1
$normal = $d.Styles.Add("Normal", 1)
